$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# Merge the two trailing runs ". " + """ into a single run ".""
# in the paragraph ending "...optional link with further information."
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute(
    "weak features (not blank, shorter than 101 characters), and an optional link with further information.",
    $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $rng1.Collapse(0)  # wdCollapseEnd
    $rng1.MoveEnd(1, 1)  # extend by one character (the closing curly quote)
    $rng1.Text = [char]0x201D
}

# --- Change 2 --------------------------------------------------------
# Replace the whole commentary paragraph about the 'work time' / workTime
# attribute with the new wording.
$oldText = "Este requisito presenta una cierta ambigüedad respecto al atributo " +
    [char]0x2018 + "work time" + [char]0x2019 + ", ya que pide solo horas, luego no deja " +
    "explícito si se debe usar un tipo de tiempo con horas, minutos y segundos o un Integer " +
    "para contar las mismas. Tras una tarea de análisis se decidió emplear un Integer ya que " +
    "aportaba la misma información, un nº de horas y presentaba una mayor sencillez respecto al código"

$newText = "En este requisito se nos pide el atributo " + [char]0x2018 + "workTime" + [char]0x2019 +
    " que obtiene un número de horas derivado de las actividades relacionadas con el mismo, " +
    "luego al ser una propiedad derivada se plantea la duda de si es necesario persistir o no " +
    "esta propiedad, tomando finalmente la decisión de no persistirla y calcularla a posteriori " +
    "en el servicio, ya que es más simple y no aporta una diferencia real notoria en el proyecto " +
    "respecto a persistirla"

$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)

$word
